$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns M ("backcrossed?") and N ("backcross_string") ------------
# Header cells: copy L1's format (style index 6 / "Neutral"-ish bold header)
# onto M1:N1, then set their text.
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1:N1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("M1").Value = "backcrossed?"
$ws.Range("N1").Value = "backcross_string"

# --- N2:N154 = "N" (the substring every row searches for in K) -----------
$ws.Range("N2:N154").Value = "N"

# --- M2:M154 = ISNUMBER(SEARCH(N,K)) --------------------------------------
# Entered in separate chunks so the saved file mirrors how this was
# actually authored (one ungrouped cell, then three fill-down batches),
# matching the shared-formula grouping boundaries.
$ws.Range("M2").Formula = "=ISNUMBER(SEARCH(N2,K2))"
$ws.Range("M3:M66").Formula = "=ISNUMBER(SEARCH(N3,K3))"
$ws.Range("M67:M130").Formula = "=ISNUMBER(SEARCH(N67,K67))"
$ws.Range("M131:M154").Formula = "=ISNUMBER(SEARCH(N131,K131))"

# --- M157 = array total of TRUE/FALSE in M2:M154 --------------------------
$ws.Range("M157").FormulaArray = "=SUM(--(M2:M154))"

# --- View state: keep K128 selected (closest reproducible approximation of
# the saved scroll/selection position recorded in the source workbook) ----
$ws.Range("K128").Select() | Out-Null
